$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CFQS")
$ws.Range("B2").Value = 400000
